$d = $word.ActiveDocument

# Locate "Group 7" in the document.
$rng = $d.Content
$found = $rng.Find.Execute("Group 7", $true, $false, $false, $false, $false,
                            $true, 1, $false, "", 0)

if ($found) {
    # $rng now spans the matched text "Group 7". Replace just the trailing
    # "7" with "C" so the text reads "Group C".
    $digit = $d.Range($rng.End - 1, $rng.End)
    $digit.Text = "C"

    # Re-touch the formatting of the replaced character so it is written out
    # as its own run (matching "Group " + "C" as two runs with identical
    # bold/underline formatting) instead of being silently merged back into
    # the preceding run.
    $digit2 = $d.Range($rng.End - 1, $rng.End)
    $digit2.Font.Bold = $false
    $digit2.Font.Bold = $true
}
